$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '72.534.93'
$ws.Range('E2').Value = '  +4.81%  '
Set-TextValue $ws.Range('D3') '4.043.08'
$ws.Range('E3').Value = '  +3.87%  '
$ws.Range('E4').Value = '  +0.11%  '
Set-TextValue $ws.Range('D5') '521.17'
$ws.Range('E5').Value = '  -0.90%  '
Set-TextValue $ws.Range('D6') '148.22'
$ws.Range('E6').Value = '  +3.83%  '
Set-TextValue $ws.Range('D7') '0.627'
$ws.Range('E7').Value = '  +2.86%  '
Set-TextValue $ws.Range('D8') '0.999'
$ws.Range('E8').Value = '  +0.01%  '
Set-TextValue $ws.Range('D9') '0.742'
$ws.Range('E9').Value = '  +3.50%  '
$ws.Range('E10').Value = '  +2.81%  '
$ws.Range('E11').Value = '  +1.27%  '
Set-TextValue $ws.Range('D12') '47.55'
$ws.Range('E12').Value = '  +13.63%  '
Set-TextValue $ws.Range('D13') '10.93'
$ws.Range('E13').Value = '  +7.56%  '
Set-TextValue $ws.Range('D14') '4.683.97'
$ws.Range('E14').Value = '  +3.66%  '
Set-TextValue $ws.Range('D15') '4.065.98'
$ws.Range('E15').Value = '  +4.54%  '
Set-TextValue $ws.Range('D16') '21.29'
$ws.Range('E16').Value = '  +8.52%  '
Set-TextValue $ws.Range('D17') '14.25'
$ws.Range('E17').Value = '  +3.33%  '
$ws.Range('E18').Value = '  -0.78%  '
$ws.Range('E19').Value = '  -1.92%  '
Set-TextValue $ws.Range('D20') '72.523.14'
$ws.Range('E20').Value = '  +4.87%  '
Set-TextValue $ws.Range('D21') '442.13'
$ws.Range('E21').Value = '  +4.28%  '
Set-TextValue $ws.Range('D22') '101.34'
$ws.Range('E22').Value = '  +15.43%  '
Set-TextValue $ws.Range('D23') '3.56'
$ws.Range('E23').Value = '  +6.80%  '
Set-TextValue $ws.Range('D24') '14.74'
$ws.Range('E24').Value = '  +4.34%  '
Set-TextValue $ws.Range('D25') '4.04'
$ws.Range('E25').Value = '  +0.06%  '
Set-TextValue $ws.Range('D26') '11.93'
$ws.Range('E26').Value = '  +3.52%  '
Set-TextValue $ws.Range('D27') '11.27'
$ws.Range('E27').Value = '  +7.33%  '
Set-TextValue $ws.Range('D28') '37.80'
$ws.Range('E28').Value = '  +4.74%  '
$ws.Range('E29').Value = '  +10.15%  '
Set-TextValue $ws.Range('D30') '13.56'
$ws.Range('E30').Value = '  +3.70%  '
Set-TextValue $ws.Range('D31') '693.98'
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('E32').Value = '  +3.61%  '
Set-TextValue $ws.Range('D33') '6.93'
$ws.Range('E33').Value = '  +17.62%  '
Set-TextValue $ws.Range('D34') '68.38'
$ws.Range('E34').Value = '  +0.47%  '
$ws.Range('E35').Value = '  +8.08%  '
Set-TextValue $ws.Range('D36') '0.442'
$ws.Range('E36').Value = '  +0.61%  '
$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D37') '41.50'
$ws.Range('E37').Value = '  +3.99%  '
$ws.Range('B38').Value = 'ThetaToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue $ws.Range('D38') '3.69'
$ws.Range('E38').Value = '  +24.93%  '
$ws.Range('E39').Value = '  +3.87%  '
Set-TextValue $ws.Range('D40') '0.999'
$ws.Range('E40').Value = '  +0.14%  '
Set-TextValue $ws.Range('D41') '0.0492'
$ws.Range('E41').Value = '  +2.77%  '
$ws.Range('E42').Value = '  -0.28%  '
$ws.Range('E43').Value = '  +5.95%  '
Set-TextValue $ws.Range('D44') '2.80'
$ws.Range('E44').Value = '  +2.02%  '
$ws.Range('E45').Value = '  +5.61%  '
$ws.Range('E46').Value = '  +5.52%  '
Set-TextValue $ws.Range('D47') '3.13'
$ws.Range('E47').Value = '  +3.39%  '
Set-TextValue $ws.Range('D48') '0.000277'
$ws.Range('E48').Value = '  +23.16%  '
Set-TextValue $ws.Range('D49') '9.13'
$ws.Range('E49').Value = '  +9.57%  '
$ws.Range('E50').Value = '  +1.50%  '
Set-TextValue $ws.Range('D51') '0.0₆0343'
$ws.Range('E51').Value = '  +0.18%  '
